# Sort the "Marks" table by the "Dept" column (ascending).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Marks")

$dataRange = $tbl.DataBodyRange
$firstRow = $dataRange.Row
$lastRow = $firstRow + $dataRange.Rows.Count - 1
$firstCol = $dataRange.Column
$lastCol = $firstCol + $dataRange.Columns.Count - 1

$deptCol = $tbl.ListColumns.Item("Dept").Range

# Table rows carry position-based formatting (e.g. the bottom border that
# belongs to the last row of the table). A plain Sort moves the formatting
# together with the cell values, which is not how Excel behaves - the
# formatting must stay anchored to the row position. Stash the "normal" row
# format and the "last row" format before sorting so they can be restored
# afterwards.
$normalFormatRow = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($firstRow, $lastCol))
$lastFormatRow = $ws.Range($ws.Cells.Item($lastRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))

$scratchRow1 = $lastRow + 50
$scratchRow2 = $lastRow + 51
$normalScratch = $ws.Range($ws.Cells.Item($scratchRow1, $firstCol), $ws.Cells.Item($scratchRow1, $lastCol))
$lastScratch = $ws.Range($ws.Cells.Item($scratchRow2, $firstCol), $ws.Cells.Item($scratchRow2, $lastCol))

$normalFormatRow.Copy()
$normalScratch.PasteSpecial(-4122)

$lastFormatRow.Copy()
$lastScratch.PasteSpecial(-4122)

# Sort the table by Dept, ascending.
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($deptCol)
$tbl.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
$tbl.Sort.Apply()

# Restore position-based row formatting.
$normalScratch.Copy()
$normalRestoreTarget = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow - 1, $lastCol))
$normalRestoreTarget.PasteSpecial(-4122)

$lastScratch.Copy()
$lastFormatRow.PasteSpecial(-4122)

# Clean up scratch cells.
$ws.Range($ws.Cells.Item($scratchRow1, $firstCol), $ws.Cells.Item($scratchRow2, $lastCol)).Clear()

$excel.CutCopyMode = $false
$ws.Range("A1").Select() | Out-Null
